$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.576331333333333
$ws.Range("H2").Value = 4.728994
$ws.Range("I2").Value = 0.01463337290888519
$ws.Range("J2").Value = 0.01463337290888519
$ws.Range("M2").Value = 0.004819666666666667
$ws.Range("N2").Value = 0.014459
$ws.Range("O2").Value = 0.04945987179224048
$ws.Range("P2").Value = 0.04945987179224049
$ws.Range("Q2").Value = 0.007597391582888888
$ws.Range("R2").Value = 0.06837652424600001
$ws.Range("S2").Value = 0.0007237647479615065
$ws.Range("T2").Value = 0.0007237647479615068

$ws.Range("G3").Value = 1.576331333333333
$ws.Range("H3").Value = 4.728994
$ws.Range("I3").Value = 0.01463337290888519
$ws.Range("J3").Value = 0.01463337290888519
$ws.Range("O3").Value = 0.5155094445470654
$ws.Range("P3").Value = 0.5155094445470654
$ws.Range("Q3").Value = 0.07918595364244445
$ws.Range("R3").Value = 0.712673582782
$ws.Range("S3").Value = 0.007543641940109477
$ws.Range("T3").Value = 0.007543641940109478

$ws.Range("G4").Value = 1.576331333333333
$ws.Range("H4").Value = 4.728994
$ws.Range("I4").Value = 0.01463337290888519
$ws.Range("J4").Value = 0.01463337290888519
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.04239200000000001
$ws.Range("N4").Value = 0.127176
$ws.Range("O4").Value = 0.4350306836606942
$ws.Range("P4").Value = 0.4350306836606942
$ws.Range("Q4").Value = 0.06682383788266667
$ws.Range("R4").Value = 0.601414540944
$ws.Range("S4").Value = 0.006365966220814204
$ws.Range("T4").Value = 0.006365966220814205

$ws.Range("G5").Value = 74.31489566666666
$ws.Range("I5").Value = 0.6898788078237544
$ws.Range("J5").Value = 0.6898788078237544
$ws.Range("M5").Value = 0.004819666666666667
$ws.Range("N5").Value = 0.014459
$ws.Range("O5").Value = 0.04945987179224048
$ws.Range("P5").Value = 0.04945987179224049
$ws.Range("Q5").Value = 0.3581730254814444
$ws.Range("R5").Value = 3.223557229333
$ws.Range("S5").Value = 0.0341213173871466
$ws.Range("T5").Value = 0.03412131738714661

$ws.Range("G6").Value = 74.31489566666666
$ws.Range("I6").Value = 0.6898788078237544
$ws.Range("J6").Value = 0.6898788078237544
$ws.Range("O6").Value = 0.5155094445470654
$ws.Range("P6").Value = 0.5155094445470654
$ws.Range("Q6").Value = 3.733159240551222
$ws.Range("S6").Value = 0.3556390410260153
$ws.Range("T6").Value = 0.3556390410260153

$ws.Range("G7").Value = 74.31489566666666
$ws.Range("I7").Value = 0.6898788078237544
$ws.Range("J7").Value = 0.6898788078237544
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.04239200000000001
$ws.Range("N7").Value = 0.127176
$ws.Range("O7").Value = 0.4350306836606942
$ws.Range("P7").Value = 0.4350306836606942
$ws.Range("Q7").Value = 3.150357057101334
$ws.Range("R7").Value = 28.353213513912
$ws.Range("S7").Value = 0.3001184494105926
$ws.Range("T7").Value = 0.3001184494105926

$ws.Range("G8").Value = 31.83044066666666
$ws.Range("H8").Value = 95.491322
$ws.Range("I8").Value = 0.2954878192673605
$ws.Range("J8").Value = 0.2954878192673605
$ws.Range("M8").Value = 0.004819666666666667
$ws.Range("N8").Value = 0.014459
$ws.Range("O8").Value = 0.04945987179224048
$ws.Range("P8").Value = 0.04945987179224049
$ws.Range("Q8").Value = 0.1534121138664444
$ws.Range("R8").Value = 1.380709024798
$ws.Range("S8").Value = 0.01461478965713238
$ws.Range("T8").Value = 0.01461478965713238

$ws.Range("G9").Value = 31.83044066666666
$ws.Range("H9").Value = 95.491322
$ws.Range("I9").Value = 0.2954878192673605
$ws.Range("J9").Value = 0.2954878192673605
$ws.Range("O9").Value = 0.5155094445470654
$ws.Range("P9").Value = 0.5155094445470654
$ws.Range("Q9").Value = 1.598980966596222
$ws.Range("R9").Value = 14.390828699366
$ws.Range("S9").Value = 0.1523267615809407
$ws.Range("T9").Value = 0.1523267615809407

$ws.Range("G10").Value = 31.83044066666666
$ws.Range("H10").Value = 95.491322
$ws.Range("I10").Value = 0.2954878192673605
$ws.Range("J10").Value = 0.2954878192673605
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.04239200000000001
$ws.Range("N10").Value = 0.127176
$ws.Range("O10").Value = 0.4350306836606942
$ws.Range("P10").Value = 0.4350306836606942
$ws.Range("Q10").Value = 1.349356040741333
$ws.Range("R10").Value = 12.144204366672
$ws.Range("S10").Value = 0.1285462680292875
$ws.Range("T10").Value = 0.1285462680292875
